$d = $word.ActiveDocument

# --- Insert the new block of paragraphs after the "(line 57)" paragraph ---
$count = $d.Paragraphs.Count
$targetIndex = 0
for ($i = 1; $i -le $count; $i++) {
  $t = $d.Paragraphs($i).Range.Text
  if ($t -eq "In ieExamples (line 57)" + [char]13) {
    $targetIndex = $i
    break
  }
}
if ($targetIndex -eq 0) {
  throw "Could not find target paragraph"
}

$targetPara = $d.Paragraphs($targetIndex).Range
$insertAt = $targetPara.End - 1
$insertPoint = $d.Range($insertAt, $insertAt)

$newParasXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>*****</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>rtPrecomputePSF</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve">DHB: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve">I put an </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>ETTBSkip</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve"> around the example (which I had just made </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>runable</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve"> rather than inline).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve">It fails </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>becuase</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve"> the routine wants to get a scene from the obscure (to me) </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>ieGetObject</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>('scene')</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>call, which doesn't work because it expects something to be set up that isn't.  The whole routine</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>thus needs a little TLC.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>*****</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>LoadRawSensorData</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve">DHB: I moved the example out of inline and made it </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>runable</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve">.  I put an </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>ETTBSkip</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve"> on it.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>It fails because it is trying</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve">to get a filename via the obscure (to me) call </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve">filename = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>vcSelectDataFile</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>('</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>stayput</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>','r'</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>);</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t xml:space="preserve">It might be that the example would run if one knew what to select here, or </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>if</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
        <w:t>the filename was just set to be something that works.  Or it might then fail below.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="15"/>
          <w:szCs w:val="15"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@

$insertPoint.InsertXML($newParasXml)

Write-Output "Inserted new paragraphs after index $targetIndex"
